$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3270934207224627
$ws.Range("C2").Value = 0.0497618483414044
$ws.Range("D2").Value = 0.0340838927683933
$ws.Range("F2").Value = 0.6015752203836726
$ws.Range("G2").Value = 0.002420357861280805
$ws.Range("K2").Value = 0.3005255237223139
$ws.Range("N2").Value = 1.26449994209834
$ws.Range("O2").Value = 2.051033743632459
$ws.Range("B3").Value = 0.2910627217096078
$ws.Range("C3").Value = 0.04581369930286883
$ws.Range("D3").Value = 0.03098075010060342
$ws.Range("F3").Value = 0.6035687897970519
$ws.Range("G3").Value = 0.002422333900540982
$ws.Range("K3").Value = 0.2632700589465173
$ws.Range("N3").Value = 1.276537187602472
$ws.Range("O3").Value = 2.068598017890963
$ws.Range("B4").Value = 0.2689476428724333
$ws.Range("C4").Value = 0.04337394687955509
$ws.Range("D4").Value = 0.02906154336075417
$ws.Range("F4").Value = 0.6051846900448723
$ws.Range("G4").Value = 0.00242361251211975
$ws.Range("K4").Value = 0.2403476157911086
$ws.Range("N4").Value = 1.284386696939677
$ws.Range("O4").Value = 2.080629194667281
$ws.Range("B5").Value = 0.2599381430966332
$ws.Range("C5").Value = 0.04237587996131253
$ws.Range("D5").Value = 0.02827600853748891
$ws.Range("F5").Value = 0.6059417242811875
$ws.Range("G5").Value = 0.002424150027636726
$ws.Range("K5").Value = 0.2309951949211069
$ws.Range("N5").Value = 1.287700812099093
$ws.Range("O5").Value = 2.085845407259228
$ws.Range("B6").Value = 0.2584422951904912
$ws.Range("C6").Value = 0.04220992130038326
$ws.Range("D6").Value = 0.02814536444959259
$ws.Range("F6").Value = 0.6060733812187138
$ws.Range("G6").Value = 0.002424240277987827
$ws.Range("K6").Value = 0.229441565045974
$ws.Range("N6").Value = 1.288258088638852
$ws.Range("O6").Value = 2.086730482618066
$ws.Range("B7").Value = 0.2688261264105165
$ws.Range("C7").Value = 0.04336050209145981
$ws.Range("D7").Value = 0.02905096324384004
$ws.Range("F7").Value = 0.6051945006677784
$ws.Range("G7").Value = 0.00242361969444852
$ws.Range("K7").Value = 0.2402215309036535
$ws.Range("N7").Value = 1.284430925020143
$ws.Range("O7").Value = 2.080698273437292
$ws.Range("B8").Value = 0.3146687320148089
$ws.Range("C8").Value = 0.04840379832420183
$ws.Range("D8").Value = 0.03301683224415086
$ws.Range("F8").Value = 0.6021812815665655
$ws.Range("G8").Value = 0.002421025671225537
$ws.Range("K8").Value = 0.2876900274418972
$ws.Range("N8").Value = 1.268555239916822
$ws.Range("O8").Value = 2.056831123455709
$ws.Range("B9").Value = 0.4046074269969608
$ws.Range("C9").Value = 0.05816757619911073
$ws.Range("D9").Value = 0.0406822876797861
$ws.Range("F9").Value = 0.5993816696541501
$ws.Range("G9").Value = 0.002416454915803331
$ws.Range("K9").Value = 0.3803776747559198
$ws.Range("N9").Value = 1.241058901525406
$ws.Range("O9").Value = 2.019923085690252
$ws.Range("B10").Value = 0.4706883870153149
$ws.Range("C10").Value = 0.06526136707520891
$ws.Range("D10").Value = 0.04624441609881558
$ws.Range("F10").Value = 0.5992216058381885
$ws.Range("G10").Value = 0.002413408421820783
$ws.Range("K10").Value = 0.4482101645114938
$ws.Range("N10").Value = 1.223069562461632
$ws.Range("O10").Value = 1.998845657667289
$ws.Range("B11").Value = 0.5007466181838538
$ws.Range("C11").Value = 0.06847066934234647
$ws.Range("D11").Value = 0.04875931925090526
$ws.Range("F11").Value = 0.5995609981261296
$ws.Range("G11").Value = 0.002412089525166103
$ws.Range("K11").Value = 0.4790069430243307
$ws.Range("N11").Value = 1.215365230350663
$ws.Range("O11").Value = 1.990569748153078
$ws.Range("B12").Value = 0.5121280383192186
$ws.Range("C12").Value = 0.06968334245539154
$ws.Range("D12").Value = 0.04970940311348215
$ws.Range("F12").Value = 0.5997488002672142
$ws.Range("G12").Value = 0.002411599675420312
$ws.Range("K12").Value = 0.4906596529526155
$ws.Range("N12").Value = 1.212516647037788
$ws.Range("O12").Value = 1.98762469656576
$ws.Range("B13").Value = 0.5096768996949379
$ws.Range("C13").Value = 0.06942228931973204
$ws.Range("D13").Value = 0.04950488665920716
$ws.Range("F13").Value = 0.5997057170516129
$ws.Range("G13").Value = 0.00241170474732922
$ws.Range("K13").Value = 0.4881504601416395
$ws.Range("N13").Value = 1.213127077181248
$ws.Range("O13").Value = 1.988250564406798
$ws.Range("B14").Value = 0.5016829977974169
$ws.Range("C14").Value = 0.06857048962575618
$ws.Range("D14").Value = 0.04883752870958347
$ws.Range("F14").Value = 0.5995752605933831
$ws.Range("G14").Value = 0.002412049032993677
$ws.Range("K14").Value = 0.4799658098476982
$ws.Range("N14").Value = 1.215129495307139
$ws.Range("O14").Value = 1.99032367120698
$ws.Range("B15").Value = 0.4967863553406744
$ws.Range("C15").Value = 0.06804839429403842
$ws.Range("D15").Value = 0.04842845749479352
$ws.Range("F15").Value = 0.599503072651089
$ws.Range("G15").Value = 0.002412261164945137
$ws.Range("K15").Value = 0.4749512366572333
$ws.Range("N15").Value = 1.216365004032106
$ws.Range("O15").Value = 1.991618107871943
$ws.Range("B16").Value = 0.4687239089488173
$ws.Range("C16").Value = 0.06505126907148906
$ws.Range("D16").Value = 0.04607974837723816
$ws.Range("F16").Value = 0.5992077210126325
$ws.Range("G16").Value = 0.002413495958771848
$ws.Range("K16").Value = 0.4461962481380226
$ws.Range("N16").Value = 1.223582698224014
$ws.Range("O16").Value = 1.999412927538117
$ws.Range("B17").Value = 0.4515074519203779
$ws.Range("C17").Value = 0.06320804338584196
$ws.Range("D17").Value = 0.04463492774293343
$ws.Range("F17").Value = 0.5991321135151679
$ws.Range("G17").Value = 0.002414270586729994
$ws.Range("K17").Value = 0.4285400182457124
$ws.Range("N17").Value = 1.228133219503157
$ws.Range("O17").Value = 2.004531017590352
$ws.Range("B18").Value = 0.4416048205511629
$ws.Range("C18").Value = 0.06214620730195008
$ws.Range("D18").Value = 0.04380246453105485
$ws.Range("F18").Value = 0.5991274245353395
$ws.Range("G18").Value = 0.002414722438255196
$ws.Range("K18").Value = 0.4183789605346817
$ws.Range("N18").Value = 1.230795655704128
$ws.Range("O18").Value = 2.007598320617049
$ws.Range("B19").Value = 0.438251948332379
$ws.Range("C19").Value = 0.06178640478371733
$ws.Range("D19").Value = 0.04352036088826594
$ws.Range("F19").Value = 0.5991325003694428
$ws.Range("G19").Value = 0.002414876511743665
$ws.Range("K19").Value = 0.414937648553547
$ws.Range("N19").Value = 1.231704856759649
$ws.Range("O19").Value = 2.008658065716602
$ws.Range("B20").Value = 0.4533401965271935
$ws.Range("C20").Value = 0.06340443039960064
$ws.Range("D20").Value = 0.0447888808234751
$ws.Range("F20").Value = 0.599136146514617
$ws.Range("G20").Value = 0.002414187473891875
$ws.Range("K20").Value = 0.4304201447089042
$ws.Range("N20").Value = 1.227644141519967
$ws.Range("O20").Value = 2.00397340397177
$ws.Range("B21").Value = 0.5040310322883386
$ws.Range("C21").Value = 0.06882075554689493
$ws.Range("D21").Value = 0.04903360952570779
$ws.Range("F21").Value = 0.5996119698814795
$ws.Range("G21").Value = 0.002411947648113946
$ws.Range("K21").Value = 0.4823700993181887
$ws.Range("N21").Value = 1.214539467289349
$ws.Range("O21").Value = 1.989709622843364
$ws.Range("B22").Value = 0.5371544140212166
$ws.Range("C22").Value = 0.07234533104512764
$ws.Range("D22").Value = 0.05179461577600364
$ws.Range("F22").Value = 0.6002684891257388
$ws.Range("G22").Value = 0.0024105396578529
$ws.Range("K22").Value = 0.5162675272928254
$ws.Range("N22").Value = 1.206376285302106
$ws.Range("O22").Value = 1.98148827938958
$ws.Range("B23").Value = 0.5194766048683732
$ws.Range("C23").Value = 0.07046562401245637
$ws.Range("D23").Value = 0.05032223630260546
$ws.Range("F23").Value = 0.5998864734350917
$ws.Range("G23").Value = 0.002411286030997318
$ws.Range("K23").Value = 0.4981810685293624
$ws.Range("N23").Value = 1.210696398133841
$ws.Range("O23").Value = 1.985775388448303
$ws.Range("B24").Value = 0.4525116271596517
$ws.Range("C24").Value = 0.06331565049420362
$ws.Range("D24").Value = 0.04471928429030925
$ws.Range("F24").Value = 0.5991342024041373
$ws.Range("G24").Value = 0.002414225029076791
$ws.Range("K24").Value = 0.4295701714314362
$ws.Range("N24").Value = 1.227865109513267
$ws.Range("O24").Value = 2.004225112218606
$ws.Range("B25").Value = 0.3802744887010192
$ws.Range("C25").Value = 0.05554000400718451
$ws.Range("D25").Value = 0.03862069642808308
$ws.Range("F25").Value = 0.5998060290066931
$ws.Range("G25").Value = 0.002417636483894417
$ws.Range("K25").Value = 0.3553481408051482
$ws.Range("N25").Value = 1.248108600443523
$ws.Range("O25").Value = 2.028847651253272
